$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 201; everything from the old row 201 downward
# (old rows 201-210) shifts down to become rows 202-211.
$ws.Rows.Item(201).Insert()

# Populate the newly inserted row 201 with its data (weekly price entry).
$ws.Cells.Item(201, 1).Value = 8
$ws.Cells.Item(201, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(201, 3).Value = "Coquimbo"
$ws.Cells.Item(201, 4).Value = 44516
$ws.Cells.Item(201, 5).Value = 4
$ws.Cells.Item(201, 6).Value = 100112032
$ws.Cells.Item(201, 7).Value = "Zapallo italiano"
$ws.Cells.Item(201, 8).Value = "Sin especificar"
$ws.Cells.Item(201, 9).Value = "Primera"
$ws.Cells.Item(201, 10).Value = 500
$ws.Cells.Item(201, 11).Value = 9500
$ws.Cells.Item(201, 12).Value = 10000
$ws.Cells.Item(201, 13).Value = 9750
$ws.Cells.Item(201, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(201, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(201, 16).Value = 139
$ws.Cells.Item(201, 17).Value = 70
$ws.Cells.Item(201, 18).Value = "Hortaliza"
